$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 "Factor" values were recalculated / re-entered.
# B2, C2, D2, F2, H2, I2 become plain (non-formula) literal values.
$ws.Range("B2").Value = 1.26
$ws.Range("C2").Value = 2.143
$ws.Range("D2").Value = 2.043
$ws.Range("F2").Value = 2.587
$ws.Range("H2").Value = 4.569
$ws.Range("I2").Value = 3.165

# E2 and G2 keep (new) formulas.
$ws.Range("E2").Formula = "=1.9"
$ws.Range("G2").Formula = "=1.287"

# Update the active selection to match the saved view state.
$ws.Range("P20").Select()
